$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new column B (log_time) - shifts page_id..key_count from B..H to C..I
$ws.Columns("B").Insert()

# Header
$ws.Range("B1").Value = "log_time"

# Row 2
$ws.Range("B2").Value = "2025-06-04 21:45:55"
$ws.Range("C2").Value = "challenge/home"
$ws.Range("D2").Value = "상세 챌린지"
$ws.Range("E2").Value = "click"
$ws.Range("F2").Value = "3일차 이런 챌린지 어때요?"
$ws.Range("G2").Value = "channel, mainTitle, click_text, list_index, challengeName, challengeSeq, list_title, chal_index, activeParticipantCount, sticker"
$ws.Range("H2").Value = "Rround, 업로드, 3일차 이런 챌린지 어때요?, 0, 업로드, Optional(104), 업로드, 0, Optional(19), RECOMMEND"
$ws.Range("I2").Value = 10

# Row 3
$ws.Range("B3").Value = "2025-06-04 21:45:55"
$ws.Range("C3").Value = "challenge/home"
$ws.Range("D3").Value = "상세 챌린지"
$ws.Range("E3").Value = "click"
$ws.Range("F3").Value = "이런 챌린지 어때요?"
$ws.Range("G3").Value = "channel, mainTitle, click_text, list_index, challengeName, challengeSeq, list_title, chal_index, activeParticipantCount, sticker"
$ws.Range("H3").Value = "Rround, 업로드, 이런 챌린지 어때요?, 0, 업로드, Optional(104), 업로드, 0, Optional(19), RECOMMEND"
$ws.Range("I3").Value = 10

# Row 4
$ws.Range("B4").Value = "2025-06-04 21:45:55"
$ws.Range("C4").Value = "challenge/challenge_detail"
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = "pageview"
$ws.Range("F4").Value = ""
$ws.Range("G4").Value = "channel, activeParticipantCount, totalFeedCount, profileSeq, challengeName "
$ws.Range("H4").Value = "Rround, Optional(19), Optional(31), Optional(627), 업로드"
$ws.Range("I4").Value = 5

# Row 5
$ws.Range("B5").Value = "2025-06-04 21:45:59"
$ws.Range("C5").Value = "challenge/challenge_detail"
$ws.Range("D5").Value = "CTA"
$ws.Range("E5").Value = "click"
$ws.Range("F5").Value = "클릭 텍스트"
$ws.Range("G5").Value = "channel, click_text, cta_text"
$ws.Range("H5").Value = "Rround, 클릭 텍스트, 인증하기"
$ws.Range("I5").Value = 3

# Row 6
$ws.Range("B6").Value = "2025-06-04 21:46:16"
$ws.Range("C6").Value = "challenge/challenge_detail"
$ws.Range("D6").Value = "상품 태그"
$ws.Range("E6").Value = "popup_click"
$ws.Range("F6").Value = "클릭 텍스트"
$ws.Range("G6").Value = "click_text, goodsId, prd_name, channel, cta_text"
$ws.Range("H6").Value = "클릭 텍스트, 2655, 660, 정관장 홍삼대정 (홍삼대정 250g * 3병), [델리스푼] 브이핏 프리미엄 이너뷰티, Rround, 상품 선택 완료"
$ws.Range("I6").Value = 5

# Row 7
$ws.Range("B7").Value = "2025-06-04 21:46:21"
$ws.Range("C7").Value = "challenge/challenge_detail"
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = "popup_imp"
$ws.Range("F7").Value = ""
$ws.Range("G7").Value = "channel, popup_title, popup_msg"
$ws.Range("H7").Value = "Rround, 인증 완료!, 다른 챌린지도 인증하고베스트 챌린저에 도전하세요 👏"
$ws.Range("I7").Value = 3

# Row 8
$ws.Range("B8").Value = "2025-06-04 21:46:25"
$ws.Range("C8").Value = "challenge/challenge_detail"
$ws.Range("D8").Value = "챌린지 상세"
$ws.Range("E8").Value = "popup_click"
$ws.Range("F8").Value = "클릭 텍스트"
$ws.Range("G8").Value = "channel, click_text, popup_msg, cta_text, popup_title"
$ws.Range("H8").Value = "Rround, 클릭 텍스트, 다른 챌린지도 인증하고베스트 챌린저에 도전하세요 👏, 확인, 인증 완료!"
$ws.Range("I8").Value = 5
